# The author swapped the contents of ppt/theme/theme1.xml (the theme driving
# the deck's slide master -- "Integral") and ppt/theme/theme2.xml (the theme
# driving the notes master -- "Office Theme"), so that the main deck now uses
# the plain "Office Theme" color scheme and the notes master keeps the old
# "Integral" colors.
#
# The font scheme (majorFont/minorFont) and format scheme (fill/line/effect/
# background styles) are byte-for-byte identical between the two themes, so
# the only real difference is the <a:clrScheme> (12 colors) -- that's what we
# reproduce here via the slide-level ThemeColorScheme, which is the
# PowerPoint object-model surface that maps onto the deck's theme part
# (ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation

function Set-ThemeRGB($colorScheme, $index, $r, $g, $b) {
    $value = $r + ($g * 256) + ($b * 65536)
    $colorScheme.Item($index).RGB = $value
}

# Target palette = the "Office Theme" clrScheme that used to live in
# ppt/theme/theme2.xml, now becoming ppt/theme/theme1.xml's palette.
$officeColors = @(
    @{ Index = 1;  Name = 'dk1';      R = 0;   G = 0;   B = 0 }
    @{ Index = 2;  Name = 'lt1';      R = 255; G = 255; B = 255 }
    @{ Index = 3;  Name = 'dk2';      R = 68;  G = 84;  B = 106 }
    @{ Index = 4;  Name = 'lt2';      R = 231; G = 230; B = 230 }
    @{ Index = 5;  Name = 'accent1';  R = 91;  G = 155; B = 213 }
    @{ Index = 6;  Name = 'accent2';  R = 237; G = 125; B = 49 }
    @{ Index = 7;  Name = 'accent3';  R = 165; G = 165; B = 165 }
    @{ Index = 8;  Name = 'accent4';  R = 255; G = 192; B = 0 }
    @{ Index = 9;  Name = 'accent5';  R = 68;  G = 114; B = 196 }
    @{ Index = 10; Name = 'accent6';  R = 112; G = 173; B = 71 }
    @{ Index = 11; Name = 'hlink';    R = 5;   G = 99;  B = 193 }
    @{ Index = 12; Name = 'folHlink'; R = 149; G = 79;  B = 114 }
)

# Apply through the first slide; the theme color scheme is shared by the
# whole deck (it is backed by the single slide-master theme part), so this
# updates ppt/theme/theme1.xml for every slide/layout/master at once.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

foreach ($entry in $officeColors) {
    Set-ThemeRGB $themeColors $entry.Index $entry.R $entry.G $entry.B
}
